# repull data, push all data, mean calculation
# Updates the dSF column (column F) values for the rows whose underlying
# data was repulled, reflecting the newly computed dSF figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = 1
    15 = -9
    16 = -2
    17 = 12
    19 = -3
    20 = -4
    22 = -1
    26 = -2
    27 = 3
    29 = -3
    31 = -4
    33 = -2
    40 = -3
    41 = -1
    42 = 1
    44 = 0
    45 = 5
    48 = -2
    51 = 1
    56 = 1
    57 = -4
    60 = 6
    63 = -1
    65 = 2
    70 = 1
    72 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
